# Edit script: swap the "B" and "C" quarter rows within each 4-row year
# block (columns A:E), and delete columns F:G (which held the now-removed
# "组合音响产销率" and "组合音响销售量" series).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs to swap (the 2nd and 3rd row of every 4-row year block).
$rowPairs = @(
    @(3,4), @(7,8), @(11,12), @(15,16),
    @(19,20), @(23,24), @(27,28), @(31,32),
    @(35,36), @(39,40), @(43,44), @(47,48),
    @(51,52), @(55,56), @(59,60), @(63,64)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rangeA = $ws.Range("A$r1" + ":E$r1")
    $rangeB = $ws.Range("A$r2" + ":E$r2")

    $valA = $rangeA.Value2
    $valB = $rangeB.Value2

    $rangeA.Value2 = $valB
    $rangeB.Value2 = $valA
}

# Remove the now-unused "组合音响产销率" (F) and "组合音响销售量" (G) columns.
$ws.Range("F1:G65").EntireColumn.Delete()
